# Daily attendance processing - 2026-01-12 17:38:25
# Swap the order of the "Recorded By" names for system-recorded sessions:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# Applies to every row in the "Recorded By" column (G) across the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
